$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 00:37:44"
$wsZhCn.Range("H4").Value = "2016-03-21 00:38:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 00:37:47"
$wsDeDe.Range("H4").Value = "2016-03-21 00:38:12"
